# Adds the "Full Series" worksheet (Year / Run, 1975-2024) after "Table 1",
# matching the Hasbrouck 2022 Kenai late-run sockeye workbook update:
#  - new sheet "Full Series" (sheetId 2) placed after "Table 1"
#  - header row: A1 "Year" (reusing the existing shared string / style),
#    B1 "Run"
#  - 50 data rows (1975-2024) with Year in col A (carrying the same cell
#    styles used on "Table 1" for the pre-1979 vs 1979+ years) and the Run
#    total in col B
#  - two trailing blank-but-styled rows (52/53), matching the source sheet
#  - selection / active-sheet view state updated to match the saved file

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- create the new sheet right after "Table 1" ------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Full Series"

# --- header row ----------------------------------------------------------
# Copy "Table 1"!A1 ("Year") formatting + shared-string value onto A1.
$ws1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122) # xlPasteFormats
$ws2.Range("A1").Value = "Year"
$ws2.Range("B1").Value = "Run"

# --- year column styles ---------------------------------------------------
# "Table 1" uses one cell style for years 1968-1978 (copy from A3, a
# non-first-in-block row so no special top-border is pulled in) and
# another for 1979 onward (copy from A13); reproduce the same two styles
# for the corresponding years here: 1975-1978 get the earlier style,
# 1979-2024 (plus the two trailing blank rows) get the later style.
$ws1.Range("A3").Copy()
$ws2.Range("A2:A5").PasteSpecial(-4122)

$ws1.Range("A13").Copy()
$ws2.Range("A6:A53").PasteSpecial(-4122)

# --- year / run data (rows 2-51 => years 1975-2024) -----------------------
$data = New-Object 'object[,]' 50,2
$data[0,0]=1975;  $data[0,1]=485350
$data[1,0]=1976;  $data[1,1]=1374607
$data[2,0]=1977;  $data[2,1]=2268567
$data[3,0]=1978;  $data[3,1]=2096342
$data[4,0]=1979;  $data[4,1]=797838
$data[5,0]=1980;  $data[5,1]=1481394
$data[6,0]=1981;  $data[6,1]=1176410
$data[7,0]=1982;  $data[7,1]=2766442
$data[8,0]=1983;  $data[8,1]=3981411
$data[9,0]=1984;  $data[9,1]=1286678
$data[10,0]=1985; $data[10,1]=2496016
$data[11,0]=1986; $data[11,1]=2945961
$data[12,0]=1987; $data[12,1]=9391896
$data[13,0]=1988; $data[13,1]=6054519
$data[14,0]=1989; $data[14,1]=6656274
$data[15,0]=1990; $data[15,1]=3224183
$data[16,0]=1991; $data[16,1]=2182082
$data[17,0]=1992; $data[17,1]=8235298
$data[18,0]=1993; $data[18,1]=4446195
$data[19,0]=1994; $data[19,1]=3886918
$data[20,0]=1995; $data[20,1]=2628555
$data[21,0]=1996; $data[21,1]=3696067
$data[22,0]=1997; $data[22,1]=4610042
$data[23,0]=1998; $data[23,1]=1902219
$data[24,0]=1999; $data[24,1]=2984568
$data[25,0]=2000; $data[25,1]=1814779
$data[26,0]=2001; $data[26,1]=2189670
$data[27,0]=2002; $data[27,1]=3466762
$data[28,0]=2003; $data[28,1]=4439571
$data[29,0]=2004; $data[29,1]=5705141
$data[30,0]=2005; $data[30,1]=6109173
$data[31,0]=2006; $data[31,1]=2848597
$data[32,0]=2007; $data[32,1]=3601777
$data[33,0]=2008; $data[33,1]=2082431
$data[34,0]=2009; $data[34,1]=2430414
$data[35,0]=2010; $data[35,1]=3596458
$data[36,0]=2011; $data[36,1]=6263091
$data[37,0]=2012; $data[37,1]=4769681
$data[38,0]=2013; $data[38,1]=3628121
$data[39,0]=2014; $data[39,1]=3404034
$data[40,0]=2015; $data[40,1]=3819016
$data[41,0]=2016; $data[41,1]=3711842
$data[42,0]=2017; $data[42,1]=2595720
$data[43,0]=2018; $data[43,1]=1867998
$data[44,0]=2019; $data[44,1]=3542442
$data[45,0]=2020; $data[45,1]=2394018
$data[46,0]=2021; $data[46,1]=3992341
$data[47,0]=2022; $data[47,1]=2929479
$data[48,0]=2023; $data[48,1]=3552933
$data[49,0]=2024; $data[49,1]=3724000

$ws2.Range("A2:B51").Value = $data

# rows 52 and 53 stay blank in column B (only the carried-over A-column
# style from the paste above applies - nothing else to set).

# --- view state: mirror the saved selections / active sheet --------------
$ws1.Range("G1:G52").Select()
$ws2.Activate()
$ws2.Range("K21").Select()
